$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns A:C slightly (closest achievable raw width to 37.85546875
# given this engine's pixel-granularity rounding on ColumnWidth).
$ws.Range("A1:C1").ColumnWidth = 37

# Row height tweaks
$ws.Range("A4:Q4").RowHeight = 28.5
$ws.Range("A5:Q5").RowHeight = 27.75
$ws.Range("A6:Q6").RowHeight = 28.5
$ws.Range("A7:Q7").RowHeight = 27
$ws.Range("A8:Q8").RowHeight = 18.75
$ws.Range("A10:Q10").RowHeight = 18.75
$ws.Range("A11:Q11").RowHeight = 17.25

# Cell value updates
$ws.Range("Q4").Value = 109
$ws.Range("Q7").Value = 12685.1
$ws.Range("P8").Value = 478225.6
$ws.Range("Q8").Value = 559503.6
$ws.Range("Q9").Value = 131.9
$ws.Range("Q10").Value = 3384.8
$ws.Range("Q11").Value = 12517.9
